$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# Update path-related values on the Settings sheet
$wsSettings.Range("B4").Value = "C:\Users\marin\AppData\Local\Programs\Python\Python38"
$wsSettings.Range("B6").Value = "C:\Users\marin\Documents\Doctorat\rpa-testing\TestingToolStable"
$wsSettings.Range("B8").Value = "C:\Users\marin\Documents\Doctorat\TestingToolStable\integrationScript.py"

# Update selection on Settings sheet to B8
$wsSettings.Activate()
$wsSettings.Range("B8").Select()

# Update row height on Constants sheet row 2
$wsConstants.Rows.Item(2).RowHeight = 30
